# Update the NPP_organised sheet: replace the formula-derived "d"-prefixed
# header cells (T1:AK1) with plain "_se"-suffixed static text labels.
$wb = $excel.ActiveWorkbook
$wsOrganised = $wb.Worksheets.Item("NPP_organised")

$cols = @("T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")
$values = @(
    "NPP_canopy_se",
    "NPP_leaf_se",
    "NPP_twigs_se",
    "NPP_flower_se",
    "NPP_fruit_se",
    "NPP_unidentified_se",
    "NPP_seed_se",
    "NPP_herbivory_se",
    "NPP_branch_turnover_se",
    "NPPacw_10cm_big_stem_se",
    "NPPacw_small_stem_se",
    "NPP_all_stem_se",
    "NPP_coarseroot_se",
    "NPP_fineroot_se",
    "NPP_herbs_se",
    "NPP_AG_se",
    "NPP_BG_se",
    "NPP_se"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsOrganised.Range($cols[$i] + "1").Value = $values[$i]
}

# Move the selection/active cell on this sheet to L27, matching the saved
# workbook view state.
$wsOrganised.Range("L27").Select()

# The "NPP" sheet's merged header cells get re-ordered (the first group of
# merges is re-created, which moves them to the end of the merge list) when
# the workbook is resaved.
$wsNpp = $wb.Worksheets.Item("NPP")
$reorderedMerges = @("K1:K2", "M1:M2", "A1:A2", "C1:C2", "E1:E2", "G1:G2", "I1:I2")
foreach ($m in $reorderedMerges) {
    $wsNpp.Range($m).UnMerge()
    $wsNpp.Range($m).Merge()
}
